$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Rows("9:30").Select()
$ws.Rows("9:30").Delete()
